$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.480.83'
$ws.Range('E2').Value = '  -6.35%  '
$ws.Range('D3').Value = '3.094.71'
$ws.Range('E3').Value = '  -6.73%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '556.05'
$ws.Range('E5').Value = '  -6.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '160.09'
$ws.Range('E6').Value = '  -12.00%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.577'
$ws.Range('E8').Value = '  -9.72%  '
$ws.Range('D9').Value = '3.091.83'
$ws.Range('E9').Value = '  -6.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.69'
$ws.Range('E10').Value = '  -2.62%  '
$ws.Range('E11').Value = '  -9.85%  '
$ws.Range('E12').Value = '  -7.31%  '
$ws.Range('D13').Value = '3.638.24'
$ws.Range('E13').Value = '  -6.53%  '
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('D15').Value = '62.542.24'
$ws.Range('E15').Value = '  -6.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '24.41'
$ws.Range('E16').Value = '  -8.86%  '
$ws.Range('D17').Value = '3.096.06'
$ws.Range('E17').Value = '  -6.81%  '
$ws.Range('E18').Value = '  -8.14%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '394.07'
$ws.Range('E19').Value = '  -8.48%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.27'
$ws.Range('E20').Value = '  -6.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.10'
$ws.Range('E21').Value = '  -7.48%  '
$ws.Range('E22').Value = '  -4.45%  '
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.67'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '67.02'
$ws.Range('E25').Value = '  -6.42%  '
$ws.Range('E26').Value = '  -5.77%  '
$ws.Range('E27').Value = '  -7.77%  '
$ws.Range('D28').Value = '0.0₃0995'
$ws.Range('E28').Value = '  -13.61%  '
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.52'
$ws.Range('E30').Value = '  -7.96%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.76'
$ws.Range('E32').Value = '  -8.91%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.74'
$ws.Range('E33').Value = '  -7.58%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.18'
$ws.Range('E34').Value = '  -6.40%  '
$ws.Range('E35').Value = '  -9.26%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '152.93'
$ws.Range('E36').Value = '  -3.95%  '
$ws.Range('E37').Value = '  -8.98%  '
$ws.Range('E38').Value = '  -9.51%  '
$ws.Range('D39').Value = '2.695.01'
$ws.Range('E39').Value = '  -6.20%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.64'
$ws.Range('E40').Value = '  -8.72%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '23.07'
$ws.Range('E41').Value = '  -12.92%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '38.06'
$ws.Range('E42').Value = '  -4.34%  '
$ws.Range('E43').Value = '  -8.82%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.686'
$ws.Range('E44').Value = '  -8.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0598'
$ws.Range('E45').Value = '  -6.95%  '
$ws.Range('E46').Value = '  -12.66%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0253'
$ws.Range('E47').Value = '  -6.83%  '
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '20.54'
$ws.Range('E49').Value = '  -10.68%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '277.31'
$ws.Range('E50').Value = '  -11.51%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0967'
$ws.Range('E51').Value = '  -5.80%  '
